$wb = $excel.ActiveWorkbook

# ===================== ALC =====================
$ws = $wb.Worksheets.Item("ALC")

# Row 40
$ws.Range("H40").Value = 7180.4546
$ws.Range("I40").Value = 4296.2
$ws.Range("J40").Value = 9584
$ws.Range("K40").Value = 4296.2
$ws.Range("L40").Value = 9584
$ws.Range("M40").Value = -4121.2
$ws.Range("N40").Value = -9934

# Row 53
$ws.Range("H53").Value = 3230.875
$ws.Range("I53").Value = 933.6667
$ws.Range("K53").Value = 933.6667
$ws.Range("M53").Value = -296.6667

# Row 64
$ws.Range("H64").Value = 18003
$ws.Range("J64").Value = 18003
$ws.Range("L64").Value = 18003
$ws.Range("N64").Value = -18499

# Row 67
$ws.Range("H67").Value = 18003
$ws.Range("J67").Value = 18003
$ws.Range("L67").Value = 18003
$ws.Range("N67").Value = -19719

# Row 70
$ws.Range("H70").Value = 64489.26
$ws.Range("I70").Value = 259749.25
$ws.Range("J70").Value = 12419.934
$ws.Range("K70").Value = 779247.75
$ws.Range("L70").Value = 37259.802
$ws.Range("M70").Value = -778977.75
$ws.Range("N70").Value = -37799.802

# Row 73
$ws.Range("H73").Value = 64489.26
$ws.Range("I73").Value = 259749.25
$ws.Range("J73").Value = 12419.934
$ws.Range("K73").Value = 779247.75
$ws.Range("L73").Value = 37259.802
$ws.Range("M73").Value = -778311.75
$ws.Range("N73").Value = -39131.802

# Row 76
$ws.Range("H76").Value = 7163
$ws.Range("I76").Value = 5847.375
$ws.Range("J76").Value = 8332.444
$ws.Range("K76").Value = 5847.375
$ws.Range("L76").Value = 8332.444
$ws.Range("M76").Value = -5532.375
$ws.Range("N76").Value = -8962.444

# Row 79
$ws.Range("H79").Value = 7163
$ws.Range("I79").Value = 5847.375
$ws.Range("J79").Value = 8332.444
$ws.Range("K79").Value = 5847.375
$ws.Range("L79").Value = 8332.444
$ws.Range("M79").Value = -4755.375
$ws.Range("N79").Value = -10516.444

# Row 111 (N111 cell removed entirely)
$ws.Range("H111").Value = 1399
$ws.Range("I111").Value = 1399
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 4197
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -1130
$ws.Range("N111").ClearContents()

# Row 116
$ws.Range("H116").Value = 7442.5
$ws.Range("I116").Value = 6198.5
$ws.Range("J116").Value = 8271.833000000001
$ws.Range("K116").Value = 6198.5
$ws.Range("L116").Value = 8271.833000000001
$ws.Range("M116").Value = -2756.5
$ws.Range("N116").Value = -15155.833

# Row 133
$ws.Range("H133").Value = 60528.8
$ws.Range("J133").Value = 60528.8
$ws.Range("L133").Value = 60528.8
$ws.Range("N133").Value = -70648.8

# ===================== ARM =====================
$ws = $wb.Worksheets.Item("ARM")

# Row 9 (N9 is new)
$ws.Range("H9").Value = 9000
$ws.Range("J9").Value = 9000
$ws.Range("L9").Value = 9000
$ws.Range("N9").Value = -9340

# Row 19
$ws.Range("H19").Value = 10103.5
$ws.Range("I19").Value = 208
$ws.Range("J19").Value = 19999
$ws.Range("K19").Value = 208
$ws.Range("L19").Value = 19999
$ws.Range("M19").Value = 21
$ws.Range("N19").Value = -20457

# Row 20 (N20 is new)
$ws.Range("H20").Value = 9000
$ws.Range("J20").Value = 9000
$ws.Range("L20").Value = 9000
$ws.Range("N20").Value = -9540

# Row 132
$ws.Range("H132").Value = 2804.2156
$ws.Range("J132").Value = 6073.778
$ws.Range("L132").Value = 18221.334
$ws.Range("N132").Value = -23281.334

# ===================== BSM =====================
$ws = $wb.Worksheets.Item("BSM")

# Row 86
$ws.Range("H86").Value = 3529
$ws.Range("I86").Value = 2885.7778
$ws.Range("J86").Value = 6423.5
$ws.Range("K86").Value = 2885.7778
$ws.Range("L86").Value = 6423.5
$ws.Range("M86").Value = -1762.7778
$ws.Range("N86").Value = -8669.5

# Row 89
$ws.Range("H89").Value = 3529
$ws.Range("I89").Value = 2885.7778
$ws.Range("J89").Value = 6423.5
$ws.Range("K89").Value = 14428.889
$ws.Range("L89").Value = 32117.5
$ws.Range("M89").Value = -8812.888999999999
$ws.Range("N89").Value = -43349.5

# Row 94
$ws.Range("H94").Value = 4754
$ws.Range("I94").Value = 4753.5
$ws.Range("K94").Value = 4753.5
$ws.Range("M94").Value = -4302.5

# Row 105
$ws.Range("H105").Value = 11147.75
$ws.Range("I105").Value = 12163.777
$ws.Range("K105").Value = 12163.777
$ws.Range("M105").Value = -10416.777

# ===================== CRP =====================
$ws = $wb.Worksheets.Item("CRP")

# Row 16
$ws.Range("H16").Value = 1348.6316
$ws.Range("I16").Value = 633.6
$ws.Range("K16").Value = 633.6
$ws.Range("M16").Value = -346.6

# Row 86
$ws.Range("H86").Value = 9302.333000000001
$ws.Range("I86").Value = 9302.333000000001
$ws.Range("K86").Value = 9302.333000000001
$ws.Range("M86").Value = -8179.333000000001

# Row 89
$ws.Range("H89").Value = 9302.333000000001
$ws.Range("I89").Value = 9302.333000000001
$ws.Range("K89").Value = 46511.665
$ws.Range("M89").Value = -40895.665

# Row 113
$ws.Range("H113").Value = 1348.6316
$ws.Range("I113").Value = 633.6
$ws.Range("K113").Value = 633.6
$ws.Range("M113").Value = 1536.4

# Row 132
$ws.Range("H132").Value = 5928.1
$ws.Range("I132").Value = 4098
$ws.Range("K132").Value = 12294
$ws.Range("M132").Value = -9764

# ===================== CUL =====================
$ws = $wb.Worksheets.Item("CUL")

# Row 4
$ws.Range("H4").Value = 11028269
$ws.Range("I4").Value = 5000239
$ws.Range("J4").Value = 23084330
$ws.Range("K4").Value = 15000717
$ws.Range("L4").Value = 69252990
$ws.Range("M4").Value = -15000605
$ws.Range("N4").Value = -69253214

# Row 131
$ws.Range("H131").Value = 7100337.5
$ws.Range("I131").Value = 25001016
$ws.Range("J131").Value = 4862752.5
$ws.Range("K131").Value = 75003048
$ws.Range("L131").Value = 14588257.5
$ws.Range("M131").Value = -74998008
$ws.Range("N131").Value = -14598337.5

# ===================== GSM =====================
$ws = $wb.Worksheets.Item("GSM")

# Row 11
$ws.Range("H11").Value = 132875624
$ws.Range("I11").Value = 6668333.5
$ws.Range("J11").Value = 208600000
$ws.Range("K11").Value = 6668333.5
$ws.Range("L11").Value = 208600000
$ws.Range("M11").Value = -6668194.5
$ws.Range("N11").Value = -208600278

# Row 113
$ws.Range("H113").Value = 4153.8887
$ws.Range("I113").Value = 2949.5
$ws.Range("K113").Value = 2949.5
$ws.Range("M113").Value = -779.5

# Row 132
$ws.Range("H132").Value = 3223.577
$ws.Range("J132").Value = 3615.923
$ws.Range("L132").Value = 10847.769
$ws.Range("N132").Value = -15907.769

# ===================== LTW =====================
$ws = $wb.Worksheets.Item("LTW")

# Row 23 (N23 cell removed entirely)
$ws.Range("H23").Value = 6000.25
$ws.Range("I23").Value = 6000.25
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 6000.25
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -5770.25
$ws.Range("N23").ClearContents()

# Row 40
$ws.Range("H40").Value = 11591.895
$ws.Range("I40").Value = 12764.8
$ws.Range("K40").Value = 12764.8
$ws.Range("M40").Value = -12628.8

# Row 46
$ws.Range("H46").Value = 6516.7095
$ws.Range("I46").Value = 5608.25
$ws.Range("J46").Value = 6832.696
$ws.Range("K46").Value = 5608.25
$ws.Range("L46").Value = 6832.696
$ws.Range("M46").Value = -5420.25
$ws.Range("N46").Value = -7208.696

# Row 68
$ws.Range("H68").Value = 4819.3
$ws.Range("I68").Value = 4099
$ws.Range("K68").Value = 4099
$ws.Range("M68").Value = -3350

# Row 71
$ws.Range("H71").Value = 4819.3
$ws.Range("I71").Value = 4099
$ws.Range("K71").Value = 20495
$ws.Range("M71").Value = -16751

# Row 100
$ws.Range("H100").Value = 11120.489
$ws.Range("I100").Value = 8692.866
$ws.Range("J100").Value = 12334.3
$ws.Range("K100").Value = 8692.866
$ws.Range("L100").Value = 12334.3
$ws.Range("M100").Value = -8151.866
$ws.Range("N100").Value = -13416.3

# ===================== WVR =====================
$ws = $wb.Worksheets.Item("WVR")

# Row 100
$ws.Range("H100").Value = 1048.7
$ws.Range("I100").Value = 871.13043
$ws.Range("J100").Value = 1632.1428
$ws.Range("K100").Value = 1742.26086
$ws.Range("L100").Value = 3264.2856
$ws.Range("M100").Value = -1201.26086
$ws.Range("N100").Value = -4346.2856
